$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Update part number and description for IC2 (M74VHC1GT04DTT1G -> M74VHC1GT14DTT1G,
# "Single inverter buffer" -> "Single Schmitt-Trigger inverter buffer")
$ws.Range("E6").Value = "M74VHC1GT14DTT1G"
$ws.Range("D6").Value = "Single Schmitt-Trigger inverter buffer"

# Update the active selection to match the author's final cursor position
$ws.Activate()
$ws.Range("E6").Select()
